$wb = $excel.ActiveWorkbook

$ws_gof = $wb.Worksheets.Item("gof")
$ws_gof.Range("D2").Value = 24393
$ws_gof.Range("F2").Value = 24459
$ws_gof.Range("G2").Value = 24645
$ws_gof.Range("D3").Value = 24323
$ws_gof.Range("F3").Value = 24469
$ws_gof.Range("G3").Value = 24882

$ws_facets = $wb.Worksheets.Item("facets")
$ws_facets.Range("B2").Value = 850
$ws_facets.Range("B3").Value = 829
$ws_facets.Range("B4").Value = 421

$ws_Estimates_1_2 = $wb.Worksheets.Item("Estimates 1-2")
$ws_Estimates_1_2.Range("B2").Value = -0.054
$ws_Estimates_1_2.Range("D2").Value = -0.046
$ws_Estimates_1_2.Range("E2").Value = 0.298
$ws_Estimates_1_2.Range("B3").Value = 0.229
$ws_Estimates_1_2.Range("D3").Value = 0.194
$ws_Estimates_1_2.Range("E3").Value = 11.015
$ws_Estimates_1_2.Range("B4").Value = 0.086
$ws_Estimates_1_2.Range("C4").Value = 0.086
$ws_Estimates_1_2.Range("D4").Value = 0.073
$ws_Estimates_1_2.Range("E4").Value = 1
$ws_Estimates_1_2.Range("B5").Value = 0.265
$ws_Estimates_1_2.Range("D5").Value = 0.224
$ws_Estimates_1_2.Range("E5").Value = 25.971
$ws_Estimates_1_2.Range("B6").Value = -0.027
$ws_Estimates_1_2.Range("D6").Value = -0.023
$ws_Estimates_1_2.Range("E6").Value = 0.057
$ws_Estimates_1_2.Range("B7").Value = 0.108
$ws_Estimates_1_2.Range("D7").Value = 0.091
$ws_Estimates_1_2.Range("E7").Value = 0.669
$ws_Estimates_1_2.Range("B8").Value = -0.284
$ws_Estimates_1_2.Range("C8").Value = 0.13
$ws_Estimates_1_2.Range("D8").Value = -0.24
$ws_Estimates_1_2.Range("E8").Value = 4.773
$ws_Estimates_1_2.Range("B9").Value = -0.081
$ws_Estimates_1_2.Range("C9").Value = 0.127
$ws_Estimates_1_2.Range("D9").Value = -0.069
$ws_Estimates_1_2.Range("E9").Value = 0.407
$ws_Estimates_1_2.Range("B10").Value = -0.058
$ws_Estimates_1_2.Range("C10").Value = 0.126
$ws_Estimates_1_2.Range("D10").Value = -0.049
$ws_Estimates_1_2.Range("E10").Value = 0.212
$ws_Estimates_1_2.Range("B11").Value = 0.151
$ws_Estimates_1_2.Range("C11").Value = 0.087
$ws_Estimates_1_2.Range("D11").Value = 0.128
$ws_Estimates_1_2.Range("E11").Value = 3.012
$ws_Estimates_1_2.Range("B12").Value = -0.01
$ws_Estimates_1_2.Range("C12").Value = 0.126
$ws_Estimates_1_2.Range("D12").Value = -0.008
$ws_Estimates_1_2.Range("E12").Value = 0.006
$ws_Estimates_1_2.Range("B13").Value = -0.056
$ws_Estimates_1_2.Range("C13").Value = 0.126
$ws_Estimates_1_2.Range("D13").Value = -0.047
$ws_Estimates_1_2.Range("E13").Value = 0.198
$ws_Estimates_1_2.Range("B14").Value = -0.225
$ws_Estimates_1_2.Range("C14").Value = 0.126
$ws_Estimates_1_2.Range("D14").Value = -0.191
$ws_Estimates_1_2.Range("E14").Value = 3.189
$ws_Estimates_1_2.Range("B15").Value = 0.231
$ws_Estimates_1_2.Range("C15").Value = 0.129
$ws_Estimates_1_2.Range("D15").Value = 0.196
$ws_Estimates_1_2.Range("E15").Value = 3.207
$ws_Estimates_1_2.Range("B16").Value = -0.19
$ws_Estimates_1_2.Range("C16").Value = 0.135
$ws_Estimates_1_2.Range("D16").Value = -0.161
$ws_Estimates_1_2.Range("E16").Value = 1.981
$ws_Estimates_1_2.Range("B17").Value = 0.397
$ws_Estimates_1_2.Range("C17").Value = 0.086
$ws_Estimates_1_2.Range("D17").Value = 0.336
$ws_Estimates_1_2.Range("E17").Value = 21.31
$ws_Estimates_1_2.Range("B18").Value = -0.351
$ws_Estimates_1_2.Range("C18").Value = 0.126
$ws_Estimates_1_2.Range("D18").Value = -0.297
$ws_Estimates_1_2.Range("E18").Value = 7.76
$ws_Estimates_1_2.Range("B19").Value = 0.158
$ws_Estimates_1_2.Range("C19").Value = 0.125
$ws_Estimates_1_2.Range("D19").Value = 0.134
$ws_Estimates_1_2.Range("E19").Value = 1.598
$ws_Estimates_1_2.Range("C20").Value = 0.079
$ws_Estimates_1_2.Range("D20").Value = 0.2
$ws_Estimates_1_2.Range("E20").Value = 8.924
$ws_Estimates_1_2.Range("B21").Value = -0.026
$ws_Estimates_1_2.Range("C21").Value = 0.131
$ws_Estimates_1_2.Range("D21").Value = -0.022
$ws_Estimates_1_2.Range("E21").Value = 0.039
$ws_Estimates_1_2.Range("B22").Value = 0.501
$ws_Estimates_1_2.Range("C22").Value = 0.506
$ws_Estimates_1_2.Range("D22").Value = 0.424
$ws_Estimates_1_2.Range("E22").Value = 0.98

$ws_Estimates_1_3 = $wb.Worksheets.Item("Estimates 1-3")
$ws_Estimates_1_3.Range("B2").Value = 0.007
$ws_Estimates_1_3.Range("C2").Value = 0.101
$ws_Estimates_1_3.Range("D2").Value = 0.006
$ws_Estimates_1_3.Range("E2").Value = 0.005
$ws_Estimates_1_3.Range("B3").Value = 0.028
$ws_Estimates_1_3.Range("C3").Value = 0.071
$ws_Estimates_1_3.Range("D3").Value = 0.024
$ws_Estimates_1_3.Range("E3").Value = 0.156
$ws_Estimates_1_3.Range("B4").Value = 0.144
$ws_Estimates_1_3.Range("C4").Value = 0.084
$ws_Estimates_1_3.Range("D4").Value = 0.122
$ws_Estimates_1_3.Range("E4").Value = 2.939
$ws_Estimates_1_3.Range("B5").Value = 0.164
$ws_Estimates_1_3.Range("C5").Value = 0.052
$ws_Estimates_1_3.Range("D5").Value = 0.139
$ws_Estimates_1_3.Range("E5").Value = 9.947
$ws_Estimates_1_3.Range("B6").Value = -0.211
$ws_Estimates_1_3.Range("C6").Value = 0.105
$ws_Estimates_1_3.Range("D6").Value = -0.179
$ws_Estimates_1_3.Range("E6").Value = 4.038
$ws_Estimates_1_3.Range("B7").Value = -0.273
$ws_Estimates_1_3.Range("C7").Value = 0.129
$ws_Estimates_1_3.Range("D7").Value = -0.231
$ws_Estimates_1_3.Range("E7").Value = 4.479
$ws_Estimates_1_3.Range("B8").Value = -0.352
$ws_Estimates_1_3.Range("C8").Value = 0.128
$ws_Estimates_1_3.Range("D8").Value = -0.298
$ws_Estimates_1_3.Range("E8").Value = 7.562
$ws_Estimates_1_3.Range("B9").Value = 0.453
$ws_Estimates_1_3.Range("C9").Value = 0.124
$ws_Estimates_1_3.Range("D9").Value = 0.384
$ws_Estimates_1_3.Range("E9").Value = 13.346
$ws_Estimates_1_3.Range("B10").Value = -0.101
$ws_Estimates_1_3.Range("C10").Value = 0.122
$ws_Estimates_1_3.Range("D10").Value = -0.086
$ws_Estimates_1_3.Range("E10").Value = 0.685
$ws_Estimates_1_3.Range("B11").Value = -0.003
$ws_Estimates_1_3.Range("C11").Value = 0.092
$ws_Estimates_1_3.Range("D11").Value = -0.003
$ws_Estimates_1_3.Range("E11").Value = 0.001
$ws_Estimates_1_3.Range("B12").Value = 0.592
$ws_Estimates_1_3.Range("C12").Value = 0.122
$ws_Estimates_1_3.Range("D12").Value = 0.501
$ws_Estimates_1_3.Range("E12").Value = 23.546
$ws_Estimates_1_3.Range("B13").Value = 0.317
$ws_Estimates_1_3.Range("C13").Value = 0.12
$ws_Estimates_1_3.Range("D13").Value = 0.268
$ws_Estimates_1_3.Range("E13").Value = 6.978
$ws_Estimates_1_3.Range("B14").Value = -0.475
$ws_Estimates_1_3.Range("C14").Value = 0.121
$ws_Estimates_1_3.Range("D14").Value = -0.402
$ws_Estimates_1_3.Range("E14").Value = 15.41
$ws_Estimates_1_3.Range("B15").Value = 0.569
$ws_Estimates_1_3.Range("C15").Value = 0.123
$ws_Estimates_1_3.Range("D15").Value = 0.482
$ws_Estimates_1_3.Range("E15").Value = 21.4
$ws_Estimates_1_3.Range("B16").Value = 0.139
$ws_Estimates_1_3.Range("C16").Value = 0.127
$ws_Estimates_1_3.Range("D16").Value = 0.118
$ws_Estimates_1_3.Range("E16").Value = 1.198
$ws_Estimates_1_3.Range("B17").Value = 0.45
$ws_Estimates_1_3.Range("C17").Value = 0.084
$ws_Estimates_1_3.Range("D17").Value = 0.381
$ws_Estimates_1_3.Range("E17").Value = 28.699
$ws_Estimates_1_3.Range("B18").Value = -0.698
$ws_Estimates_1_3.Range("C18").Value = 0.118
$ws_Estimates_1_3.Range("D18").Value = -0.591
$ws_Estimates_1_3.Range("E18").Value = 34.99
$ws_Estimates_1_3.Range("I18").Value = 0.999
$ws_Estimates_1_3.Range("B19").Value = 0.194
$ws_Estimates_1_3.Range("C19").Value = 0.118
$ws_Estimates_1_3.Range("D19").Value = 0.164
$ws_Estimates_1_3.Range("E19").Value = 2.703
$ws_Estimates_1_3.Range("B20").Value = -0.064
$ws_Estimates_1_3.Range("C20").Value = 0.078
$ws_Estimates_1_3.Range("D20").Value = -0.054
$ws_Estimates_1_3.Range("E20").Value = 0.673
$ws_Estimates_1_3.Range("B21").Value = -0.158
$ws_Estimates_1_3.Range("C21").Value = 0.122
$ws_Estimates_1_3.Range("D21").Value = -0.134
$ws_Estimates_1_3.Range("E21").Value = 1.677
$ws_Estimates_1_3.Range("B22").Value = 0.724
$ws_Estimates_1_3.Range("C22").Value = 0.489
$ws_Estimates_1_3.Range("D22").Value = 0.613
$ws_Estimates_1_3.Range("E22").Value = 2.192

$ws_Estimates_2_3 = $wb.Worksheets.Item("Estimates 2-3")
$ws_Estimates_2_3.Range("B2").Value = 0.061
$ws_Estimates_2_3.Range("D2").Value = 0.052
$ws_Estimates_2_3.Range("E2").Value = 0.387
$ws_Estimates_2_3.Range("B3").Value = -0.201
$ws_Estimates_2_3.Range("D3").Value = -0.17
$ws_Estimates_2_3.Range("E3").Value = 9.275
$ws_Estimates_2_3.Range("B4").Value = 0.058
$ws_Estimates_2_3.Range("C4").Value = 0.088
$ws_Estimates_2_3.Range("D4").Value = 0.049
$ws_Estimates_2_3.Range("E4").Value = 0.434
$ws_Estimates_2_3.Range("B5").Value = -0.101
$ws_Estimates_2_3.Range("D5").Value = -0.086
$ws_Estimates_2_3.Range("E5").Value = 3.632
$ws_Estimates_2_3.Range("B6").Value = -0.184
$ws_Estimates_2_3.Range("C6").Value = 0.121
$ws_Estimates_2_3.Range("D6").Value = -0.156
$ws_Estimates_2_3.Range("E6").Value = 2.312
$ws_Estimates_2_3.Range("B7").Value = -0.38
$ws_Estimates_2_3.Range("C7").Value = 0.135
$ws_Estimates_2_3.Range("D7").Value = -0.322
$ws_Estimates_2_3.Range("E7").Value = 7.923
$ws_Estimates_2_3.Range("B8").Value = -0.067
$ws_Estimates_2_3.Range("D8").Value = -0.057
$ws_Estimates_2_3.Range("E8").Value = 0.258
$ws_Estimates_2_3.Range("B9").Value = 0.534
$ws_Estimates_2_3.Range("D9").Value = 0.452
$ws_Estimates_2_3.Range("E9").Value = 16.873
$ws_Estimates_2_3.Range("B10").Value = -0.042
$ws_Estimates_2_3.Range("D10").Value = -0.036
$ws_Estimates_2_3.Range("E10").Value = 0.106
$ws_Estimates_2_3.Range("B11").Value = -0.154
$ws_Estimates_2_3.Range("D11").Value = -0.13
$ws_Estimates_2_3.Range("E11").Value = 3.443
$ws_Estimates_2_3.Range("B12").Value = 0.602
$ws_Estimates_2_3.Range("C12").Value = 0.131
$ws_Estimates_2_3.Range("D12").Value = 0.51
$ws_Estimates_2_3.Range("E12").Value = 21.118
$ws_Estimates_2_3.Range("B13").Value = 0.373
$ws_Estimates_2_3.Range("C13").Value = 0.131
$ws_Estimates_2_3.Range("D13").Value = 0.316
$ws_Estimates_2_3.Range("E13").Value = 8.107
$ws_Estimates_2_3.Range("B14").Value = -0.25
$ws_Estimates_2_3.Range("C14").Value = 0.131
$ws_Estimates_2_3.Range("D14").Value = -0.212
$ws_Estimates_2_3.Range("E14").Value = 3.642
$ws_Estimates_2_3.Range("B15").Value = 0.339
$ws_Estimates_2_3.Range("C15").Value = 0.135
$ws_Estimates_2_3.Range("D15").Value = 0.287
$ws_Estimates_2_3.Range("E15").Value = 6.306
$ws_Estimates_2_3.Range("B16").Value = 0.329
$ws_Estimates_2_3.Range("C16").Value = 0.143
$ws_Estimates_2_3.Range("D16").Value = 0.279
$ws_Estimates_2_3.Range("E16").Value = 5.293
$ws_Estimates_2_3.Range("B17").Value = 0.052
$ws_Estimates_2_3.Range("C17").Value = 0.088
$ws_Estimates_2_3.Range("D17").Value = 0.044
$ws_Estimates_2_3.Range("E17").Value = 0.349
$ws_Estimates_2_3.Range("B18").Value = -0.347
$ws_Estimates_2_3.Range("C18").Value = 0.133
$ws_Estimates_2_3.Range("D18").Value = -0.294
$ws_Estimates_2_3.Range("E18").Value = 6.807
$ws_Estimates_2_3.Range("B19").Value = 0.036
$ws_Estimates_2_3.Range("C19").Value = 0.131
$ws_Estimates_2_3.Range("D19").Value = 0.03
$ws_Estimates_2_3.Range("E19").Value = 0.076
$ws_Estimates_2_3.Range("B20").Value = -0.301
$ws_Estimates_2_3.Range("C20").Value = 0.08
$ws_Estimates_2_3.Range("D20").Value = -0.255
$ws_Estimates_2_3.Range("E20").Value = 14.156
$ws_Estimates_2_3.Range("B21").Value = -0.133
$ws_Estimates_2_3.Range("C21").Value = 0.14
$ws_Estimates_2_3.Range("D21").Value = -0.113
$ws_Estimates_2_3.Range("E21").Value = 0.902
$ws_Estimates_2_3.Range("B22").Value = 0.223
$ws_Estimates_2_3.Range("C22").Value = 0.523
$ws_Estimates_2_3.Range("D22").Value = 0.189
$ws_Estimates_2_3.Range("E22").Value = 0.182

$ws_Main_effect_1_2 = $wb.Worksheets.Item("Main effect 1-2")
$ws_Main_effect_1_2.Range("B2").Value = 0.688
$ws_Main_effect_1_2.Range("C2").Value = 0.582
$ws_Main_effect_1_2.Range("B3").Value = 0.49
$ws_Main_effect_1_2.Range("C3").Value = 0.415

$ws_Main_effect_1_3 = $wb.Worksheets.Item("Main effect 1-3")
$ws_Main_effect_1_3.Range("B2").Value = 0.295
$ws_Main_effect_1_3.Range("C2").Value = 0.25
$ws_Main_effect_1_3.Range("B3").Value = 0.214
$ws_Main_effect_1_3.Range("C3").Value = 0.181

$ws_Main_effect_2_3 = $wb.Worksheets.Item("Main effect 2-3")
$ws_Main_effect_2_3.Range("B2").Value = -0.392
$ws_Main_effect_2_3.Range("C2").Value = -0.332
$ws_Main_effect_2_3.Range("B3").Value = -0.276
$ws_Main_effect_2_3.Range("C3").Value = -0.234
